$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.080.45'
$ws.Range('E2').Value = '  -3.04%  '

# Row 3
$ws.Range('D3').Value = '2.376.49'
$ws.Range('E3').Value = '  -3.65%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '546.05'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.25%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.79'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.57%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('E8').Value = '  -0.69%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.153'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.61%  '

# Row 10
$ws.Range('E10').Value = '  -1.38%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.322'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.84%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.70'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.54%  '

# Row 13
$ws.Range('D13').Value = '67.006.21'
$ws.Range('E13').Value = '  -2.96%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000165'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.73%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.52'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.66%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '10.15'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -5.78%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '323.96'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -5.58%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.68'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -5.43%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.70'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.62%  '

# Row 20
$ws.Range('E20').Value = '  -0.42%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.81'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.63%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.14'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.94%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.53'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.17%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.81'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.75%  '

# Row 25
$ws.Range('D25').Value = '0.0₃0780'
$ws.Range('E25').Value = '  -4.74%  '

# Row 26
$ws.Range('E26').Value = '  -4.14%  '

# Row 27
$ws.Range('E27').Value = '  +0.06%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '408.19'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -7.15%  '

# Row 29
$ws.Range('E29').Value = '  -3.65%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.56'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.28%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.85'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.37%  '

# Row 32
$ws.Range('E32').Value = '  -0.74%  '

# Row 33
$ws.Range('E33').Value = '  -0.08%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.46'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.49%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.102'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.94%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.288'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.74%  '

# Row 37
$ws.Range('E37').Value = '  -3.33%  '

# Row 38
$ws.Range('E38').Value = '  -7.01%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.04'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.41%  '

# Row 40
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '126.92'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.67%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.23'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.87%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.92'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -7.73%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0701'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.49%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.464'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.13%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.545'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.07%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0904'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.43%  '

# Row 47
$ws.Range('E47').Value = '  -1.24%  '

# Row 48
$ws.Range('E48').Value = '  -9.00%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.17'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.47%  '

# Row 50
$ws.Range('B50').Value = 'Hedera'
$ws.Range('C50').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0420'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.51%  '

# Row 51
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0198'
$ws.Range('E51').Value = '  -7.22%  '
